# py-polars empty.xlsx fixture update (#16802)
#
# The workbook used to contain a single, totally empty sheet ("Sheet1").
# The fixture is reshaped into two sheets so the io tests can cover both
# "no data at all" and "header row but zero data rows" cases:
#   - "no_data": the original empty sheet, just renamed.
#   - "no_rows": a new sheet with a bold header row (colx, coly, colz)
#                and no data rows underneath it.

$wb = $excel.ActiveWorkbook

# Rename the existing (empty) sheet to "no_data".
$noData = $wb.Worksheets.Item(1)
$noData.Name = "no_data"

# Add the new "no_rows" sheet right after "no_data" and give it a header-only row.
$noRows = $wb.Worksheets.Add($null, $noData)
$noRows.Name = "no_rows"

$noRows.Range("A1").Value = "colx"
$noRows.Range("B1").Value = "coly"
$noRows.Range("C1").Value = "colz"
$noRows.Range("A1:C1").Font.Bold = $true

# Keep the original "no_data" sheet as the active/selected one, matching the
# unmodified state of the first worksheet part in the source workbook.
$noData.Activate()
